# This script applies a weekly data-refresh to the "Perejil" (Vega Monumental
# Concepción) sheet: every existing weekly record (rows 130..221, laid out as
# two rows per record: "Primera" / "Segunda" quality) is shifted down by one
# record (2 rows), the newest record is written into rows 130-131 with new
# figures, and the oldest record (which was in rows 220-221) is appended again
# at the bottom (rows 222-223) as the sheet grows by one record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: capture the current ("old") values of every column that varies
# record to record, for every row in the block that will shift (130..221).
# ---------------------------------------------------------------------------
$cols = @("D","J","K","L","M","P","O")
$oldVals = @{}
foreach ($c in $cols) {
    $oldVals[$c] = @{}
    for ($r = 130; $r -le 221; $r++) {
        $oldVals[$c][$r] = $ws.Range("$c$r").Value2
    }
}

# ---------------------------------------------------------------------------
# Step 2: append two new rows (222, 223) at the bottom of the table by
# duplicating the current last record (rows 220-221) cell by cell, including
# all the columns that never change between records (A, B, C, E, F, G, H, I,
# N, Q, R). The D column also needs its date number format copied across so
# the new cells render the same way as the rest of the date column.
# ---------------------------------------------------------------------------
$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($c in $allCols) {
    $ws.Range("$c" + "222").Value = $ws.Range("$c" + "220").Value2
    $ws.Range("$c" + "223").Value = $ws.Range("$c" + "221").Value2
}
$ws.Range("D222").NumberFormat = $ws.Range("D220").NumberFormat
$ws.Range("D223").NumberFormat = $ws.Range("D221").NumberFormat

# ---------------------------------------------------------------------------
# Step 3: shift every record from rows 132..221 down by one record (2 rows)
# into rows 134..223, i.e. new row n gets the value previously held by row
# (n - 2). Processed so every destination is only written once, using the
# values captured in step 1 (so write order does not matter).
# ---------------------------------------------------------------------------
for ($r = 223; $r -ge 132; $r--) {
    $src = $r - 2
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $oldVals[$c][$src]
    }
}

# ---------------------------------------------------------------------------
# Step 4: write the brand-new weekly record into rows 130-131 (this week's
# "Primera" / "Segunda" quality prices for Perejil).
# ---------------------------------------------------------------------------
$ws.Range("D130").Value = 45062
$ws.Range("J130").Value = 150
$ws.Range("K130").Value = 700
$ws.Range("L130").Value = 800
$ws.Range("M130").Value = 767
$ws.Range("P130").Value = 767
$ws.Range("O130").Value = $oldVals["O"][130]

$ws.Range("D131").Value = 45062
$ws.Range("J131").Value = 100
$ws.Range("K131").Value = 600
$ws.Range("L131").Value = 600
$ws.Range("M131").Value = 600
$ws.Range("P131").Value = 600
$ws.Range("O131").Value = $oldVals["O"][131]
